$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.4555622440282

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "60bd88b8fc436774352f53b9"
$ws.Range("E3").Value = "Annes"
$ws.Range("G3").Value = 13.49322534775249

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "60b1742bce2b39e0f1d19a1a"
$ws.Range("E4").Value = "Sabrina"
$ws.Range("G4").Value = 13.0399847539705

$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "5f2c1a97a6809c060fec8820"
$ws.Range("E5").Value = "Maggie"
$ws.Range("G5").Value = 13.02496522416173

$ws.Range("G6").Value = 12.06802414159128

$ws.Range("C7").Value = 7
$ws.Range("D7").Value = "5ff3974450a7199965624df7"
$ws.Range("E7").Value = "Anh"
$ws.Range("G7").Value = 10.28276308106538

$ws.Range("C8").Value = 6
$ws.Range("D8").Value = "60b7cd4be2d4cc6bb252e016"
$ws.Range("E8").Value = "Chris"
$ws.Range("G8").Value = 10.06553779822961

$ws.Range("G9").Value = 8.109626341792918

$ws.Range("G10").Value = 5.328503579847605

$ws.Range("C11").Value = 15
$ws.Range("D11").Value = "60b76ee2219ac1ce25ccea43"
$ws.Range("E11").Value = "Richie"
$ws.Range("F11").Value = "male"
$ws.Range("G11").Value = 2.335090669544204

$ws.Range("C12").Value = 14
$ws.Range("D12").Value = "60186dc2cc1aa8103499603a"
$ws.Range("E12").Value = "Emily"
$ws.Range("F12").Value = "female"
$ws.Range("G12").Value = 2.285979158260564

$ws.Range("G13").Value = 1.402438053655847

$ws.Range("G14").Value = 15.4707565882647

$ws.Range("G15").Value = 11.24716022619973

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = "60778ed0fde3e9c3a96f1d11"
$ws.Range("E16").Value = "Melissa"
$ws.Range("F16").Value = "female"
$ws.Range("G16").Value = 10.28396858577988

$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "60ba8ba51a5e0a105396888a"
$ws.Range("E17").Value = "Alfredo"
$ws.Range("F17").Value = "male"
$ws.Range("G17").Value = 10.17922288895593

$ws.Range("G18").Value = 9.11564123155263

$ws.Range("G19").Value = 7.090080371342193

$ws.Range("C20").Value = 7
$ws.Range("D20").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("E20").Value = "Katherine"
$ws.Range("G20").Value = 5.492612310008454

$ws.Range("C21").Value = 8
$ws.Range("D21").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("E21").Value = "Valeria"
$ws.Range("G21").Value = 5.467819815978743

$ws.Range("G22").Value = 5.255963789825259

$ws.Range("G23").Value = 4.443507100185125

$ws.Range("G24").Value = 3.499128519851419

$ws.Range("G25").Value = 3.153071154622522

